$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 62, pushing rows 62:87 down to 63:88.
$ws.Rows.Item(62).Insert()

# Copy formatting from the row above (row 61) isn't required for values; explicitly
# set the new row's cell values/format to match the template used by this dataset.
$ws.Range("A62").Value = 9
$ws.Range("B62").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C62").Value = "Metropolitana"
$ws.Range("D62").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44636)
$ws.Range("D62").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = 100112005
$ws.Range("G62").Value = "Puerro"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 97
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = 7000
$ws.Range("N62").Value = "`$/paquete 20 unidades"
$ws.Range("O62").Value = "Provincia de Chacabuco"
$ws.Range("P62").Value = 350
$ws.Range("Q62").Value = 20
$ws.Range("R62").Value = "Hortaliza"
